# Updated with Deduplication of Names Formula
#
# - Adds a "dedup names" CSE array-formula helper in column J (J2:J12),
#   pulling unique names out of L10:L14.
# - Extends the existing P-column "Extract" CSE array formula (previously a
#   shared formula referencing $L$10:$L$15) down through row 27, now
#   referencing $L$10:$L$16 (rows 1-24) and finally a broken-reference
#   tail (rows 25-27) left over from the column having been dragged past
#   its helper range.
# - Grows Table1 and the hidden _FilterDatabase name by one row
#   (L9:N15 -> L9:N16).
# - Widens column L and moves the active selection to J10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Resize Table1 (L9:N15 -> L9:N16) -----------------------------------
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("L9:N16"))

# --- Update the hidden _xlnm._FilterDatabase name -----------------------
foreach ($n in $ws.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$L`$10:`$L`$16"
    }
}

# --- Widen column L -------------------------------------------------------
$ws.Columns.Item(12).ColumnWidth = 15.8776041666667

# --- Column J: dedup-names helper formulas (J2:J12) -----------------------
# J2 pulls from the (slightly mis-sized) $L$10:$L$19 block; J3:J12 reference
# the tighter $L$10:$L$14 block used for the rest of the fill.
$ws.Range("J2").FormulaArray = '=IFERROR(INDEX($L$10:$L$19, MATCH(0, COUNTIF($J$1:J1, $L$10:$L$19), 0)), "")'

for ($r = 3; $r -le 12; $r++) {
    $prev = $r - 1
    $formula = '=IFERROR(INDEX($L$10:$L$14, MATCH(0, COUNTIF($J$1:J' + $prev + ', $L$10:$L$14), 0)), "")'
    $ws.Range("J$r").FormulaArray = $formula
}

# --- Column P: extend the "Extract" array formula down to row 24 ----------
# Rows 1-15 reference ROW(L<r>); rows 16-24 reference ROW(L<r+1>) (the
# off-by-one left behind when the fill handle was dragged past the bottom
# of the lookup range).
for ($r = 1; $r -le 15; $r++) {
    $formula = '=IFERROR(INDEX($L$10:$L$16,SMALL(IF(ISTEXT($L$10:$L$16),ROW($L$1:$L$14)), ROW(L' + $r + '))),"")'
    $ws.Range("P$r").FormulaArray = $formula
}

for ($r = 16; $r -le 24; $r++) {
    $target = $r + 1
    $formula = '=IFERROR(INDEX($L$10:$L$16,SMALL(IF(ISTEXT($L$10:$L$16),ROW($L$1:$L$14)), ROW(L' + $target + '))),"")'
    $ws.Range("P$r").FormulaArray = $formula
}

# --- Column P rows 25-27: trailing broken-reference formulas --------------
for ($r = 25; $r -le 27; $r++) {
    $formula = '=IFERROR(INDEX($L$10:$L$19,SMALL(IF(ISTEXT($L$10:$L$19),ROW($L$10:$L$19)), ROW(#REF!))),"")'
    $ws.Range("P$r").Formula = $formula
}

# --- Active selection moves to J10 -----------------------------------------
$ws.Range("J10").Select()
